$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report header text (volume/number and week-covering dates)
$ws.Cells.Item(8, 1).Value = "Volume 30   Number  45"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

function Set-TextCellFromStyle($row, $col, $text, $styleRow, $styleCol) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $src = $ws.Cells.Item($styleRow, $styleCol)
    $src.Copy()
    $c.PasteSpecial(-4122) | Out-Null
}

function Set-NumCellFromStyle($row, $col, $value, $styleRow, $styleCol) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $value
    $src = $ws.Cells.Item($styleRow, $styleCol)
    $src.Copy()
    $c.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

Set-TextCellFromStyle 14 3 "0" 14 1
Set-TextCellFromStyle 14 7 "0" 14 1
Set-TextCellFromStyle 14 8 "***.*" 14 1
$ws.Cells.Item(14, 14).Value2 = -76.190476190476
Set-TextCellFromStyle 15 6 "0" 14 1
$ws.Cells.Item(15, 7).Value2 = 2
$ws.Cells.Item(15, 8).Value2 = -100
$ws.Cells.Item(15, 12).Value2 = 12.903225806451
$ws.Cells.Item(15, 13).Value2 = 25
$ws.Cells.Item(15, 14).Value2 = -57.317073170731
$ws.Cells.Item(16, 3).Value2 = 4
$ws.Cells.Item(16, 4).Value2 = 18
$ws.Cells.Item(16, 5).Value2 = -77.777777777777
$ws.Cells.Item(16, 6).Value2 = 41
$ws.Cells.Item(16, 7).Value2 = 64
$ws.Cells.Item(16, 8).Value2 = -35.9375
$ws.Cells.Item(16, 9).Value2 = 476
$ws.Cells.Item(16, 10).Value2 = 623
$ws.Cells.Item(16, 11).Value2 = -23.595505617977
$ws.Cells.Item(16, 12).Value2 = 8.675799086757
$ws.Cells.Item(16, 13).Value2 = 11.475409836065
$ws.Cells.Item(16, 14).Value2 = -73.031161473087
$ws.Cells.Item(17, 3).Value2 = 27
$ws.Cells.Item(17, 4).Value2 = 20
$ws.Cells.Item(17, 5).Value2 = 35
$ws.Cells.Item(17, 6).Value2 = 82
$ws.Cells.Item(17, 7).Value2 = 82
$ws.Cells.Item(17, 8).Value2 = 0
$ws.Cells.Item(17, 9).Value2 = 958
$ws.Cells.Item(17, 10).Value2 = 928
$ws.Cells.Item(17, 11).Value2 = 3.232758620689
$ws.Cells.Item(17, 12).Value2 = 25.392670157068
$ws.Cells.Item(17, 13).Value2 = 98.343685300207
$ws.Cells.Item(17, 14).Value2 = -20.365752285951
$ws.Cells.Item(18, 3).Value2 = 6
$ws.Cells.Item(18, 5).Value2 = -33.333333333333
$ws.Cells.Item(18, 6).Value2 = 24
$ws.Cells.Item(18, 7).Value2 = 45
$ws.Cells.Item(18, 8).Value2 = -46.666666666666
$ws.Cells.Item(18, 9).Value2 = 397
$ws.Cells.Item(18, 10).Value2 = 380
$ws.Cells.Item(18, 11).Value2 = 4.473684210526
$ws.Cells.Item(18, 12).Value2 = 66.108786610878
$ws.Cells.Item(18, 13).Value2 = 71.120689655172
$ws.Cells.Item(18, 14).Value2 = -78.162816281628
$ws.Cells.Item(19, 3).Value2 = 21
$ws.Cells.Item(19, 4).Value2 = 12
$ws.Cells.Item(19, 5).Value2 = 75
$ws.Cells.Item(19, 6).Value2 = 67
$ws.Cells.Item(19, 7).Value2 = 63
$ws.Cells.Item(19, 8).Value2 = 6.349206349206
$ws.Cells.Item(19, 9).Value2 = 698
$ws.Cells.Item(19, 10).Value2 = 834
$ws.Cells.Item(19, 11).Value2 = -16.306954436450
$ws.Cells.Item(19, 12).Value2 = 6.240487062404
$ws.Cells.Item(19, 13).Value2 = 110.240963855422
$ws.Cells.Item(19, 14).Value2 = -2.649930264993
$ws.Cells.Item(20, 3).Value2 = 9
$ws.Cells.Item(20, 4).Value2 = 8
$ws.Cells.Item(20, 5).Value2 = 12.5
$ws.Cells.Item(20, 6).Value2 = 35
$ws.Cells.Item(20, 7).Value2 = 24
$ws.Cells.Item(20, 8).Value2 = 45.833333333333
$ws.Cells.Item(20, 9).Value2 = 415
$ws.Cells.Item(20, 10).Value2 = 333
$ws.Cells.Item(20, 11).Value2 = 24.624624624624
$ws.Cells.Item(20, 12).Value2 = 192.253521126761
$ws.Cells.Item(20, 13).Value2 = 126.775956284153
$ws.Cells.Item(20, 14).Value2 = -63.787085514834
$ws.Cells.Item(21, 3).Value2 = 67
$ws.Cells.Item(21, 4).Value2 = 67
$ws.Cells.Item(21, 5).Value2 = 0
$ws.Cells.Item(21, 6).Value2 = 251
$ws.Cells.Item(21, 7).Value2 = 280
$ws.Cells.Item(21, 8).Value2 = -10.357142857142
$ws.Cells.Item(21, 9).Value2 = 2994
$ws.Cells.Item(21, 10).Value2 = 3151
$ws.Cells.Item(21, 11).Value2 = -4.982545223738
$ws.Cells.Item(21, 12).Value2 = 30.856643356643
$ws.Cells.Item(21, 13).Value2 = 76.325088339222
$ws.Cells.Item(21, 14).Value2 = -55.931704445098
$ws.Cells.Item(22, 3).Value2 = 1
$ws.Cells.Item(22, 5).Value2 = 0
$ws.Cells.Item(22, 6).Value2 = 6
$ws.Cells.Item(22, 7).Value2 = 2
$ws.Cells.Item(22, 8).Value2 = 200
$ws.Cells.Item(22, 9).Value2 = 34
$ws.Cells.Item(22, 10).Value2 = 34
$ws.Cells.Item(22, 12).Value2 = 78.947368421052
$ws.Cells.Item(22, 13).Value2 = -8.108108108108
$ws.Cells.Item(23, 3).Value2 = 3
$ws.Cells.Item(23, 4).Value2 = 2
$ws.Cells.Item(23, 5).Value2 = 50
$ws.Cells.Item(23, 6).Value2 = 6
$ws.Cells.Item(23, 7).Value2 = 6
$ws.Cells.Item(23, 8).Value2 = 0
$ws.Cells.Item(23, 9).Value2 = 71
$ws.Cells.Item(23, 10).Value2 = 69
$ws.Cells.Item(23, 11).Value2 = 2.898550724637
$ws.Cells.Item(23, 12).Value2 = 0
$ws.Cells.Item(23, 13).Value2 = 77.5
$ws.Cells.Item(24, 3).Value2 = 31
$ws.Cells.Item(24, 4).Value2 = 46
$ws.Cells.Item(24, 5).Value2 = -32.608695652173
$ws.Cells.Item(24, 6).Value2 = 132
$ws.Cells.Item(24, 7).Value2 = 166
$ws.Cells.Item(24, 8).Value2 = -20.481927710843
$ws.Cells.Item(24, 9).Value2 = 1802
$ws.Cells.Item(24, 10).Value2 = 1694
$ws.Cells.Item(24, 11).Value2 = 6.375442739079
$ws.Cells.Item(24, 12).Value2 = 51.047778709136
$ws.Cells.Item(24, 13).Value2 = 49.792186201163
$ws.Cells.Item(25, 3).Value2 = 21
$ws.Cells.Item(25, 4).Value2 = 29
$ws.Cells.Item(25, 5).Value2 = -27.586206896551
$ws.Cells.Item(25, 7).Value2 = 78
$ws.Cells.Item(25, 8).Value2 = 21.794871794871
$ws.Cells.Item(25, 9).Value2 = 1082
$ws.Cells.Item(25, 10).Value2 = 1116
$ws.Cells.Item(25, 11).Value2 = -3.046594982078
$ws.Cells.Item(25, 12).Value2 = 20.355951056729
$ws.Cells.Item(25, 13).Value2 = -1.457194899817
Set-NumCellFromStyle 26 4 1 14 6
Set-NumCellFromStyle 26 5 -100 14 11
$ws.Cells.Item(26, 6).Value2 = 3
$ws.Cells.Item(26, 8).Value2 = 0
$ws.Cells.Item(26, 10).Value2 = 53
$ws.Cells.Item(26, 11).Value2 = 15.094339622641
$ws.Cells.Item(26, 12).Value2 = 24.489795918367
$ws.Cells.Item(27, 3).Value2 = 2
$ws.Cells.Item(27, 4).Value2 = 2
$ws.Cells.Item(27, 9).Value2 = 116
$ws.Cells.Item(27, 10).Value2 = 104
$ws.Cells.Item(27, 11).Value2 = 11.538461538461
$ws.Cells.Item(27, 12).Value2 = 31.818181818181
Set-TextCellFromStyle 28 4 "0" 14 1
Set-TextCellFromStyle 28 5 "***.*" 14 1
Set-TextCellFromStyle 28 6 "0" 14 1
$ws.Cells.Item(28, 8).Value2 = -100
$ws.Cells.Item(28, 12).Value2 = -33.333333333333
$ws.Cells.Item(28, 13).Value2 = 35.135135135135
$ws.Cells.Item(28, 14).Value2 = -69.135802469135
Set-TextCellFromStyle 29 4 "0" 14 1
Set-TextCellFromStyle 29 5 "***.*" 14 1
Set-TextCellFromStyle 29 6 "0" 14 1
$ws.Cells.Item(29, 8).Value2 = -100
$ws.Cells.Item(29, 12).Value2 = 33.333333333333
$ws.Cells.Item(29, 13).Value2 = 33.333333333333
$ws.Cells.Item(29, 14).Value2 = -72.602739726027
$ws.Cells.Item(30, 12).Value2 = -66.666666666666

$excel.CutCopyMode = 0

